$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 5
    4  = 3
    5  = 5
    6  = 8
    7  = 6
    8  = 6
    9  = 6
    10 = 4
    11 = 7
    12 = 9
    13 = 6
    14 = 6
    15 = 6
    16 = 5
    17 = 5
    18 = 7
    19 = 3
    20 = 3
    21 = 4
    22 = 12
    23 = 3
    24 = 9
    25 = 4
    26 = 4
    27 = 6
    28 = 5
    29 = 2
    30 = 2
    31 = 3
    32 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
